$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.499.12"
$ws.Range("E2").Value = "  +4.50%  "

# Row 3
$ws.Range("D3").Value = "1.599.70"
$ws.Range("E3").Value = "  +2.17%  "

# Row 4
$ws.Range("E4").Value = "  -0.31%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.11"
$ws.Range("E5").Value = "  +2.24%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.499"
$ws.Range("E6").Value = "  +1.74%  "

# Row 7
$ws.Range("E7").Value = "  -0.32%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.19"
$ws.Range("E8").Value = "  +9.76%  "

# Row 9
$ws.Range("E9").Value = "  +1.27%  "

# Row 10
$ws.Range("E10").Value = "  +0.95%  "

# Row 11
$ws.Range("E11").Value = "  +2.20%  "

# Row 12
$ws.Range("D12").Value = "1.825.51"
$ws.Range("E12").Value = "  +2.01%  "

# Row 13
$ws.Range("D13").Value = "1.598.77"
$ws.Range("E13").Value = "  +2.26%  "

# Row 14
$ws.Range("E14").Value = "  +0.95%  "

# Row 15
$ws.Range("E15").Value = "  +3.31%  "

# Row 16
$ws.Range("D16").Value = "28.491.56"
$ws.Range("E16").Value = "  +4.50%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.44"
$ws.Range("E17").Value = "  +2.50%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.40"
$ws.Range("E18").Value = "  +7.49%  "

# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.55"
$ws.Range("E19").Value = "  +1.33%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0713"
$ws.Range("E20").Value = "  +1.58%  "

# Row 21
$ws.Range("E21").Value = "  -0.22%  "

# Row 22
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.44"
$ws.Range("E23").Value = "  +2.33%  "

# Row 24
$ws.Range("E24").Value = "  +1.14%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.83"
$ws.Range("E25").Value = "  -0.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.33"
$ws.Range("E26").Value = "  +2.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.64"
$ws.Range("E27").Value = "  +0.00%  "

# Row 28
$ws.Range("E28").Value = "  +1.33%  "

# Row 29
$ws.Range("E29").Value = "  -0.37%  "

# Row 30
$ws.Range("E30").Value = "  +0.97%  "

# Row 31
$ws.Range("E31").Value = "  +1.13%  "

# Row 32
$ws.Range("E32").Value = "  +0.55%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("E33").Value = "  +1.03%  "

# Row 34
$ws.Range("D34").Value = "1.424.22"
$ws.Range("E34").Value = "  -0.82%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  -0.23%  "

# Row 36
$ws.Range("E36").Value = "  -4.06%  "

# Row 37
$ws.Range("E37").Value = "  -0.22%  "

# Row 38
$ws.Range("E38").Value = "  +0.93%  "

# Row 39
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.546"
$ws.Range("E39").Value = "  +2.41%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.53"
$ws.Range("E40").Value = "  +8.18%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.824"
$ws.Range("E41").Value = "  +2.08%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.77"
$ws.Range("E42").Value = "  -2.82%  "

# Row 44
$ws.Range("E44").Value = "  +6.36%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.981"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.04"
$ws.Range("E46").Value = "  +0.83%  "

# Row 47
$ws.Range("D47").Value = "1.737.12"
$ws.Range("E47").Value = "  +1.98%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.66"
$ws.Range("E48").Value = "  +1.92%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.13"
$ws.Range("E49").Value = "  -0.20%  "

# Row 50
$ws.Range("E50").Value = "  +5.40%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0527"
$ws.Range("E51").Value = "  +0.32%  "

